# Updates the cryptos price list (columns D = Price, E = Volume(1h)) row by row.
# Values must stay plain TEXT cells (matching the sheet's original inlineStr
# cells) rather than being auto-coerced to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # A bare numeric-looking string (e.g. "1.003") gets auto-converted to a
    # Number by Value assignment, which would change the stored cell type.
    # Prefixing with a single quote forces Excel to keep it as text (exactly
    # like typing it into the grid); re-applying the "Normal" style afterwards
    # clears the quote-prefix formatting flag Excel sets along the way, so the
    # cell style is left exactly as it was.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '28.810.22'
$ws.Range("E2").Value = '  +2.52%  '
$ws.Range("D3").Value = '1.875.45'
$ws.Range("E3").Value = '  +2.05%  '
Set-TextValue $ws.Range("D4") '1.003'
$ws.Range("E4").Value = '  +0.21%  '
Set-TextValue $ws.Range("D5") '324.85'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +0.29%  '
Set-TextValue $ws.Range("D7") '0.4611'
$ws.Range("E7").Value = '  -0.67%  '
Set-TextValue $ws.Range("D8") '0.3865'
$ws.Range("E8").Value = '  -0.28%  '
Set-TextValue $ws.Range("D9") '0.07876'
$ws.Range("E9").Value = '  +0.03%  '
Set-TextValue $ws.Range("D10") '0.9834'
$ws.Range("E10").Value = '  +2.10%  '
Set-TextValue $ws.Range("D11") '21.83'
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").Value = '1.891.67'
Set-TextValue $ws.Range("D13") '7.002'
$ws.Range("E13").Value = '  +1.11%  '
Set-TextValue $ws.Range("D14") '5.706'
$ws.Range("E14").Value = '  +0.20%  '
Set-TextValue $ws.Range("D15") '0.06969'
$ws.Range("E15").Value = '  +2.21%  '
Set-TextValue $ws.Range("D16") '88.47'
$ws.Range("E16").Value = '  +0.35%  '
Set-TextValue $ws.Range("D17") '1.002'
$ws.Range("E17").Value = '  +0.17%  '
Set-TextValue $ws.Range("D18") '0.00001004'
$ws.Range("E18").Value = '  +0.79%  '
Set-TextValue $ws.Range("D19") '16.79'
$ws.Range("E19").Value = '  +0.35%  '
Set-TextValue $ws.Range("D20") '1.005'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").Value = '28.791.81'
$ws.Range("E21").Value = '  +2.39%  '
Set-TextValue $ws.Range("D22") '5.277'
$ws.Range("E22").Value = '  -0.95%  '
Set-TextValue $ws.Range("D23") '11.09'
$ws.Range("E23").Value = '  +0.56%  '
Set-TextValue $ws.Range("D24") '2.102'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '2.116.20'
$ws.Range("E25").Value = '  +2.77%  '
Set-TextValue $ws.Range("D26") '153.04'
$ws.Range("E26").Value = '  -1.14%  '
Set-TextValue $ws.Range("D27") '19.31'
$ws.Range("E27").Value = '  +0.38%  '
Set-TextValue $ws.Range("D28") '5.850'
$ws.Range("E28").Value = '  +2.88%  '
Set-TextValue $ws.Range("D29") '1.992'
$ws.Range("E29").Value = '  +1.33%  '
Set-TextValue $ws.Range("D30") '119.13'
$ws.Range("E30").Value = '  +0.70%  '
Set-TextValue $ws.Range("D31") '0.09331'
$ws.Range("E31").Value = '  +0.85%  '
Set-TextValue $ws.Range("D32") '0.9226'
$ws.Range("E32").Value = '  -1.64%  '
Set-TextValue $ws.Range("D33") '5.303'
$ws.Range("E33").Value = '  +0.48%  '
Set-TextValue $ws.Range("D34") '1.338'
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("E35").Value = '  +0.47%  '
Set-TextValue $ws.Range("D36") '0.05784'
$ws.Range("E36").Value = '  -1.66%  '
Set-TextValue $ws.Range("D37") '1.149'
Set-TextValue $ws.Range("D38") '0.02081'
$ws.Range("E38").Value = '  -2.44%  '
Set-TextValue $ws.Range("D39") '7.664'
$ws.Range("E39").Value = '  -1.55%  '
Set-TextValue $ws.Range("D40") '0.5642'
$ws.Range("E40").Value = '  +0.61%  '
Set-TextValue $ws.Range("D41") '0.1783'
$ws.Range("E41").Value = '  +1.04%  '
Set-TextValue $ws.Range("D42") '9.785'
$ws.Range("E42").Value = '  -1.28%  '
Set-TextValue $ws.Range("D43") '0.07211'
$ws.Range("E43").Value = '  -0.64%  '
Set-TextValue $ws.Range("D44") '11.76'
$ws.Range("E44").Value = '  +0.54%  '
Set-TextValue $ws.Range("D45") '0.5305'
$ws.Range("E45").Value = '  +0.38%  '
Set-TextValue $ws.Range("D46") '2.148'
$ws.Range("E46").Value = '  +0.49%  '
Set-TextValue $ws.Range("D47") '1.122'
$ws.Range("E47").Value = '  -0.99%  '
Set-TextValue $ws.Range("D48") '1.838'
$ws.Range("E48").Value = '  +0.35%  '
Set-TextValue $ws.Range("D49") '113.38'
$ws.Range("E49").Value = '  +0.43%  '
Set-TextValue $ws.Range("D50") '2.414'
$ws.Range("E50").Value = '  +3.73%  '
$ws.Range("E51").Value = '  +0.29%  '
